$wb = $excel.ActiveWorkbook

# --- Sheet "Overview": update Latest HO Xliff Generate Date for rows 4 and 5 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-16 10:15:40"
$wsOverview.Range("G5").Value = "2016-08-16 10:15:40"

# --- Sheet "zh-cn": update Priority (rows 4 & 5) and Correspond Handoff/Handback Datetime (rows 4 & 5) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H4").Value = "2016-08-16 10:15:34"
$wsZhCn.Range("H5").Value = "2016-08-16 10:15:34"
$wsZhCn.Range("K4").Value = "2016-08-16 10:15:52"
$wsZhCn.Range("K5").Value = "2016-08-16 10:15:52"

# --- Sheet "de-de": update Priority (rows 4 & 5), Correspond Handoff Datetime (rows 4 & 5),
#     and Correspond Handback Datetime (rows 4 & 5) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H4").Value = "2016-08-16 10:15:40"
$wsDeDe.Range("H5").Value = "2016-08-16 10:15:40"
$wsDeDe.Range("K4").Value = "2016-08-16 10:16:01"
$wsDeDe.Range("K5").Value = "2016-08-16 10:16:01"
